$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mark "Y" in columns C (Handles processor ${Property}?) and D (Accumulates
# messages in For() loop?) for the commands that were updated to support
# these features.
$rows = @(78, 107, 136, 170, 198, 206, 228)

foreach ($r in $rows) {
    $ws.Range("C$r").Value = "Y"
    $ws.Range("D$r").Value = "Y"
}

# Update the active cell/selection to the last-edited cell.
$ws.Range("D199").Select()

# Update page setup: fit-to-page printing, landscape orientation, scaled to 62%.
$ws.PageSetup.Zoom = 62
$ws.PageSetup.FitToPagesTall = 0
$ws.PageSetup.Orientation = 2
